$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format while writing, so values like
# "571.41" are stored as literal strings instead of being auto-parsed
# into numbers by the COM layer (matches the source inlineStr cells).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '65.546.73'
$ws.Range("E2").Value = '  -0.26%  '
$ws.Range("D3").Value = '3.172.67'
$ws.Range("E3").Value = '  -5.18%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '571.41'
$ws.Range("E5").Value = '  -1.35%  '
$ws.Range("D6").Value = '171.35'
$ws.Range("E6").Value = '  -3.62%  '
$ws.Range("D8").Value = '0.600'
$ws.Range("E8").Value = '  -3.35%  '
$ws.Range("D9").Value = '3.169.37'
$ws.Range("E9").Value = '  -5.23%  '
$ws.Range("E10").Value = '  -4.00%  '
$ws.Range("D11").Value = '6.55'
$ws.Range("E11").Value = '  -4.92%  '
$ws.Range("E12").Value = '  -5.04%  '
$ws.Range("D13").Value = '3.712.34'
$ws.Range("E13").Value = '  -5.55%  '
$ws.Range("D14").Value = '0.136'
$ws.Range("E14").Value = '  +0.97%  '
$ws.Range("D15").Value = '27.26'
$ws.Range("E15").Value = '  -5.24%  '
$ws.Range("D16").Value = '65.490.58'
$ws.Range("E16").Value = '  -0.39%  '
$ws.Range("E17").Value = '  -3.91%  '
$ws.Range("D18").Value = '3.161.28'
$ws.Range("E18").Value = '  -5.11%  '
$ws.Range("D19").Value = '5.71'
$ws.Range("E19").Value = '  -0.50%  '
$ws.Range("D20").Value = '12.85'
$ws.Range("E20").Value = '  -4.60%  '
$ws.Range("D21").Value = '358.36'
$ws.Range("E21").Value = '  -1.51%  '
$ws.Range("D22").Value = '7.27'
$ws.Range("E22").Value = '  -2.75%  '
$ws.Range("E23").Value = '  +0.39%  '
$ws.Range("D24").Value = '69.08'
$ws.Range("E24").Value = '  -3.80%  '
$ws.Range("D25").Value = '0.494'
$ws.Range("E25").Value = '  -5.52%  '
$ws.Range("D26").Value = '3.288.42'
$ws.Range("E26").Value = '  -6.52%  '
$ws.Range("D27").Value = '0.0000115'
$ws.Range("E27").Value = '  -6.36%  '
$ws.Range("D28").Value = '9.79'
$ws.Range("E28").Value = '  +1.57%  '
$ws.Range("D29").Value = '0.177'
$ws.Range("E29").Value = '  -1.38%  '
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.30%  '
$ws.Range("E32").Value = '  -2.26%  '
$ws.Range("D33").Value = '5.36'
$ws.Range("E33").Value = '  -5.75%  '
$ws.Range("D34").Value = '21.95'
$ws.Range("E34").Value = '  -3.99%  '
$ws.Range("B35").Value = 'Fetch.AI'
$ws.Range("C35").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D35").Value = '1.20'
$ws.Range("E35").Value = '  -2.07%  '
$ws.Range("B36").Value = 'Aptos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D36").Value = '6.59'
$ws.Range("E36").Value = '  -4.50%  '
$ws.Range("D37").Value = '159.40'
$ws.Range("E37").Value = '  -0.26%  '
$ws.Range("D38").Value = '1.45'
$ws.Range("E38").Value = '  -4.01%  '
$ws.Range("D39").Value = '0.833'
$ws.Range("E39").Value = '  -2.16%  '
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").Value = '1.79'
$ws.Range("E40").Value = '  +1.44%  '
$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").Value = '26.39'
$ws.Range("E41").Value = '  -2.95%  '
$ws.Range("D42").Value = '2.49'
$ws.Range("E42").Value = '  -2.98%  '
$ws.Range("D43").Value = '2.642.91'
$ws.Range("E43").Value = '  -0.92%  '
$ws.Range("E44").Value = '  -2.84%  '
$ws.Range("D45").Value = '4.19'
$ws.Range("E45").Value = '  -2.44%  '
$ws.Range("D46").Value = '39.65'
$ws.Range("E46").Value = '  -0.25%  '
$ws.Range("D47").Value = '0.0659'
$ws.Range("E47").Value = '  -1.21%  '
$ws.Range("D48").Value = '24.09'
$ws.Range("E48").Value = '  -1.32%  '
$ws.Range("D49").Value = '327.62'
$ws.Range("E49").Value = '  -2.76%  '
$ws.Range("D50").Value = '0.0274'
$ws.Range("E50").Value = '  -2.28%  '
$ws.Range("E51").Value = '  -1.56%  '

# Restore the original (default) cell formatting on column D so we don't
# leave a stray text-format style behind that wasn't in the source file.
$ws.Range("D2:D51").ClearFormats()

